$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '90.506.18'
$ws.Range("E2").Value = '  -0.93%  '

# Row 3
$ws.Range("D3").Value = '3.134.24'
$ws.Range("E3").Value = '  +0.70%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.03'
$ws.Range("E5").Value = '  +7.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '644.26'
$ws.Range("E6").Value = '  +3.43%  '

# Row 7
$ws.Range("E7").Value = '  +11.34%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.359'
$ws.Range("E8").Value = '  -5.87%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("D10").Value = '3.131.91'
$ws.Range("E10").Value = '  +0.67%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.726'
$ws.Range("E11").Value = '  +1.36%  '

# Row 12
$ws.Range("E12").Value = '  +4.24%  '

# Row 13
$ws.Range("E13").Value = '  +5.69%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").Value = '  -3.96%  '

# Row 15
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.63'
$ws.Range("E15").Value = '  +4.47%  '

# Row 16
$ws.Range("D16").Value = '90.172.05'
$ws.Range("E16").Value = '  -1.11%  '

# Row 17
$ws.Range("D17").Value = '3.708.02'
$ws.Range("E17").Value = '  +0.27%  '

# Row 18
$ws.Range("D18").Value = '3.118.56'
$ws.Range("E18").Value = '  -0.15%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.67'
$ws.Range("E19").Value = '  -1.41%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.47'
$ws.Range("E20").Value = '  +2.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000213'
$ws.Range("E21").Value = '  -3.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '450.96'
$ws.Range("E22").Value = '  +3.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  +10.87%  '

# Row 24
$ws.Range("E24").Value = '  +3.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.03'
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.43'
$ws.Range("E26").Value = '  +4.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.40'
$ws.Range("E27").Value = '  +1.76%  '

# Row 28
$ws.Range("E28").Value = '  +0.42%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.88'
$ws.Range("E30").Value = '  +8.62%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.161'
$ws.Range("E31").Value = '  -3.95%  '

# Row 32
$ws.Range("E32").Value = '  +34.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.38'
$ws.Range("E33").Value = '  +15.21%  '

# Row 34
$ws.Range("E34").Value = '  +2.64%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.151'
$ws.Range("E35").Value = '  +6.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '514.28'
$ws.Range("E36").Value = '  -2.32%  '

# Row 37
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("E37").Value = '  +5.76%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.12'
$ws.Range("E38").Value = '  +0.41%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.32'
$ws.Range("E39").Value = '  +2.99%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.424'
$ws.Range("E40").Value = '  +11.84%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0868'
$ws.Range("E41").Value = '  +0.40%  '

# Row 42
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.20'
$ws.Range("E42").Value = '  -0.46%  '

# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("B44").Value = 'Binance-PegBSC-USD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.744'
$ws.Range("E44").Value = '  -16.68%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.34'
$ws.Range("E45").Value = '  +42.19%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.94'
$ws.Range("E46").Value = '  +1.71%  '

# Row 47
$ws.Range("E47").Value = '  +13.70%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '149.34'
$ws.Range("E48").Value = '  +1.74%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.59'
$ws.Range("E49").Value = '  +9.74%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.24'
$ws.Range("E50").Value = '  +2.58%  '

# Row 51
$ws.Range("E51").Value = '  +4.06%  '
